$d = $word.ActiveDocument

# Step 1: remove the leading empty paragraph before "8. Mai"
$p108 = $d.Paragraphs.Item(108)
$p108.Range.Delete()

# Step 2: insert a new empty paragraph right after the "Wie in der..." paragraph
$p109 = $d.Paragraphs.Item(109)
$insPoint1 = $d.Range($p109.Range.End, $p109.Range.End)
$emptyXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint1.InsertXML($emptyXml)

# Step 3: insert the new diary content after that new empty paragraph
$p110 = $d.Paragraphs.Item(110)
$insPoint2 = $d.Range($p110.Range.End, $p110.Range.End)
$mainXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>09. - 10. Mai</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Die letzten Features werden dem Spiel zugefügt: eine Hintergrundmusik, Töne zu den Button- und Card Clicks sowie Cheat codes und der Bestrafung derer. Die MP3SPI Library zum Abspielen vom mp3-Dateien bereitet zuerst einige Schwierigkeiten, welche jedoch mit der Hilfe von Forumsbeiträgen erfolgreich überwunden werden können. Ausserdem wird der Fall behandelt, wenn ein Spieler ein Game verlässt, ohne jedoch die App zu verlassen, da bisher ein komplettes Logout die einzige Möglichkeit zum Verlassen eines Games war. Auch auf GUI-Seite müssen noch gewisse Änderungen angebracht werden, zum Beispiel zum Muten der Musik und Soundeffekte. Ausserdem wird der Entwurf für den Trailer gemacht.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>11. - 13. Mai</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Auch wenn ab jetzt nur noch „der letzte Schliff“ zu tun ist (und die Präsentation), braucht dies mehr Zeit als gedacht. Wir spielen das Spiel mehrere Male, und immer wieder tauchen Fehler auf - speziell beim Verlassen eines Games tauchen immer wieder Fälle auf, welche wir nicht beachtet hatten. Da beim Verlassen des Games das ArrayList mit den verbleibenden Spielern verschoben wird, braucht man je nach dem eine Anpassung der Variable, welche angibt, wer am Zug ist. Ausserdem muss natürlich der Zug beendet werden, falls der Spieler, der das Game verlässt, selbst gerade am Zug war. Speziell durch das ursprüngliche Design der Methoden, welche das Verlassen des Games verarbeiten – welche wir auf zwei Klassen aufgeteilt hatten, weil das uns zu Beginn am L</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">ogischsten </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">schien – werden Anpassungen erschwert. Schlussendlich schaffen wir es, alle Operationen in der richtigen Reihenfolge und mit Berücksichtigung aller Spezialfälle durchzuführen: die eine Methode, welche den Spieler aus dem Spiel nimmt, gibt ein boolean Array an die andere Methode zurück, welche sich darum kümmert, alle anderen Spieler darüber zu informieren. </w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Obwohl das Beheben dieser kleinen Fehler nur sehr wenig Codearbeit erfordert, stellt sich schnell heraus, dass Veränderungen in einem für uns bereits so komplexen Code nur mit äusserster Sorgfalt durchgeführt können. Nur so können wir gleichzeitig kollaterale Schäden verhindern und alle Codeteile finden, welche angepasst werden müssen. Auch die Cheat codes brauchen einige Anpassungen – von einer eigenen Methode in der Game Klasse zu einem neuen Feld in der Player Klasse, welche angibt, ob der Spieler bereits einen Cheat benutzt hat.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>Gameplay wird aufgenommen und der Trailer wird erstellt. Der QA Report wird zum Schluss im Plenum diskutiert und gemeinsam abgeschlossen. Dadurch kann jeder seine Sicht zu den QA Massnahmen einbringen und seine Meinung dazu äussern, was uns die Messungen zeigen. Parallel werden die Slides der PowerPoint abgeschlossen und anschliessend die Präsentation noch einige Male geübt. Ein Problem mit der JAR Datei taucht noch auf, welches nicht richtig auf die Ressourcen zugreift: Sound und Icons werden nicht geladen. In der Übungsstunde wird auch noch dieses Problem behoben.</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insPoint2.InsertXML($mainXml)

# Step 4: re-create the _GoBack bookmark at the split point between "am L" and "ogischsten"
$found = $d.Content
$found.Find.ClearFormatting()
$found.Find.Execute("ogischsten", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($found.Start, $found.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
